# Sync automático del tracker - 2025-10-01 18:47:09 UTC
#
# 1) Rows 197-199: mark as Completed with results, outcome, profit, ROI and
#    the send timestamp.
# 2) Rows 200-212: append newly scraped "UEFA Europa League" fixtures with
#    Status = Pending (result columns left blank, matching new predictions).
# 3) Sheet dimension grows from A1:Q199 to A1:Q212 automatically once the
#    new cells are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Resolve pending predictions for rows 197-199
# ---------------------------------------------------------------------------
$completed = @(
    @{ Row = 197; Result = "Draw";     Outcome = "Fallo";   Profit = -3.2; ROI = -100; Sent = "2025-10-01 04:26:50" },
    @{ Row = 198; Result = "Home Win"; Outcome = "Acierto"; Profit = 1.6;  ROI = 80;   Sent = "2025-10-01 04:26:50" },
    @{ Row = 199; Result = "Away Win"; Outcome = "Acierto"; Profit = 0.75; ROI = 125;  Sent = "2025-10-01 04:26:50" }
)

foreach ($entry in $completed) {
    $r = $entry.Row
    $ws.Cells.Item($r, 12).Value = "Completed"        # L - Status
    $ws.Cells.Item($r, 13).Value = $entry.Result      # M - Result
    $ws.Cells.Item($r, 14).Value = $entry.Outcome     # N - Resultado_Real
    $ws.Cells.Item($r, 15).Value = $entry.Profit       # O - Profit
    $ws.Cells.Item($r, 16).Value = $entry.ROI         # P - ROI
    $ws.Cells.Item($r, 17).Value = $entry.Sent        # Q - Enviado
}

# ---------------------------------------------------------------------------
# 2) Append new pending predictions (rows 200-212)
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 200; Date = "2025-10-02"; Liga = "UEFA Europa League"; Local = "Celtic";             Visitante = "SC Braga";            Prediccion = "Home Win"; Probabilidad = "65.66%"; Cuota = 1.91; EV = "24.15%"; Stake = 1.8; StakePct = 0.02791515253207479; Kelly = 0.2791515253207478 },
    @{ Row = 201; Date = "2025-10-02"; Liga = "UEFA Europa League"; Local = "Panathinaikos";      Visitante = "GO Ahead Eagles";     Prediccion = "Home Win"; Probabilidad = "81.52%"; Cuota = 1.57; EV = "26.70%"; Stake = 3.1; StakePct = 0.04909559739067899; Kelly = 0.4909559739067899 },
    @{ Row = 202; Date = "2025-10-02"; Liga = "UEFA Europa League"; Local = "Ludogorets";         Visitante = "Real Betis";          Prediccion = "Away Win"; Probabilidad = "70.52%"; Cuota = 1.83; EV = "27.76%"; Stake = 2.2; StakePct = 0.03500238967882782; Kelly = 0.3500238967882781 },
    @{ Row = 203; Date = "2025-10-02"; Liga = "UEFA Europa League"; Local = "Bologna";            Visitante = "SC Freiburg";         Prediccion = "Home Win"; Probabilidad = "66.01%"; Cuota = 1.9;  EV = "24.16%"; Stake = 1.8; StakePct = 0.02824007439207552; Kelly = 0.2824007439207552 },
    @{ Row = 204; Date = "2025-10-02"; Liga = "UEFA Europa League"; Local = "Fenerbahce";         Visitante = "Nice";                Prediccion = "Home Win"; Probabilidad = "73.21%"; Cuota = 1.73; EV = "25.39%"; Stake = 2.3; StakePct = 0.03651991158815074; Kelly = 0.3651991158815074 },
    @{ Row = 205; Date = "2025-10-02"; Liga = "UEFA Europa League"; Local = "Plzen";              Visitante = "Malmo FF";            Prediccion = "Home Win"; Probabilidad = "64.00%"; Cuota = 1.95; EV = "23.56%"; Stake = 1.6; StakePct = 0.02610924401215796; Kelly = 0.2610924401215796 },
    @{ Row = 206; Date = "2025-10-02"; Liga = "UEFA Europa League"; Local = "AS Roma";            Visitante = "Lille";               Prediccion = "Home Win"; Probabilidad = "74.97%"; Cuota = 1.7;  EV = "26.18%"; Stake = 2.5; StakePct = 0.03922194584395017; Kelly = 0.3922194584395016 },
    @{ Row = 207; Date = "2025-10-02"; Liga = "UEFA Europa League"; Local = "FC Basel 1893";      Visitante = "VfB Stuttgart";       Prediccion = "Away Win"; Probabilidad = "67.10%"; Cuota = 1.83; EV = "21.57%"; Stake = 1.7; StakePct = 0.0274689875466225;  Kelly = 0.274689875466225  },
    @{ Row = 208; Date = "2025-10-02"; Liga = "UEFA Europa League"; Local = "Genk";               Visitante = "Ferencvarosi TC";     Prediccion = "Home Win"; Probabilidad = "76.61%"; Cuota = 1.7;  EV = "28.94%"; Stake = 2.7; StakePct = 0.04319713003957059; Kelly = 0.4319713003957059 },
    @{ Row = 209; Date = "2025-10-02"; Liga = "UEFA Europa League"; Local = "Lyon";               Visitante = "Red Bull Salzburg";   Prediccion = "Home Win"; Probabilidad = "83.41%"; Cuota = 1.55; EV = "27.99%"; Stake = 3.2; StakePct = 0.05;                 Kelly = 0.5323940662311937 },
    @{ Row = 210; Date = "2025-10-02"; Liga = "UEFA Europa League"; Local = "FC Porto";           Visitante = "FK Crvena Zvezda";    Prediccion = "Home Win"; Probabilidad = "90.59%"; Cuota = 1.4;  EV = "25.56%"; Stake = 3.2; StakePct = 0.05;                 Kelly = 0.6706630702058514 },
    @{ Row = 211; Date = "2025-10-02"; Liga = "UEFA Europa League"; Local = "Nottingham Forest";  Visitante = "FC Midtjylland";      Prediccion = "Home Win"; Probabilidad = "87.44%"; Cuota = 1.48; EV = "28.12%"; Stake = 3.2; StakePct = 0.05;                 Kelly = 0.6128013083179003 },
    @{ Row = 212; Date = "2025-10-02"; Liga = "UEFA Europa League"; Local = "Celta Vigo";         Visitante = "PAOK";                Prediccion = "Home Win"; Probabilidad = "73.63%"; Cuota = 1.73; EV = "26.11%"; Stake = 2.4; StakePct = 0.0375062794028851;  Kelly = 0.375062794028851  }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    # Date/percent-looking text must stay literal text (as in the rest of
    # the sheet), so force Text format before assigning - otherwise Excel's
    # autodetection would turn them into a date serial / a numeric percent.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value  = $entry.Date          # A  - Date
    $ws.Cells.Item($r, 2).Value  = $entry.Liga          # B  - Liga
    $ws.Cells.Item($r, 3).Value  = $entry.Local         # C  - Local
    $ws.Cells.Item($r, 4).Value  = $entry.Visitante     # D  - Visitante
    $ws.Cells.Item($r, 5).Value  = $entry.Prediccion    # E  - Prediccion
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value  = $entry.Probabilidad  # F  - Probabilidad
    $ws.Cells.Item($r, 7).Value  = $entry.Cuota         # G  - Cuota_Bet365
    $ws.Cells.Item($r, 8).NumberFormat = "@"
    $ws.Cells.Item($r, 8).Value  = $entry.EV            # H  - EV
    $ws.Cells.Item($r, 9).Value  = $entry.Stake         # I  - Stake
    $ws.Cells.Item($r, 10).Value = $entry.StakePct      # J  - StakePct
    $ws.Cells.Item($r, 11).Value = $entry.Kelly         # K  - KellyFrac
    $ws.Cells.Item($r, 12).Value = "Pending"            # L  - Status
    # M (Result), N (Resultado_Real), O (Profit), P (ROI), Q (Enviado)
    # stay blank for newly-added pending predictions.
}
